$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correción error reporte ventas
# Reorder the "Materias primas" ingredient lists so the last ingredient
# (Vainilla / Limon) is listed first, matching the corrected report order.

$ws.Range("C2").Value = "1.0-Vainilla (ml),1.0-Huevos (unidad),3.0-Leche (litros),2.0-Harina  (kg),"
$ws.Range("C4").Value = "1.0-Vainilla (ml),2.0-Huevos (unidad),5.0-Harina  (kg),"
$ws.Range("C6").Value = "2.0-Limon (unidad),5.0-Huevos (unidad),5.0-Crema (litros),4.0-Harina  (kg),"
$ws.Range("C7").Value = "0.1-Vainilla (ml),2.0-Huevos (unidad),0.2-Leche (litros),0.3-Harina  (kg),"
